# Apply "Minor Updates" changes to Multibeam_metadata workbook:
#  - Replace old C-SCAMP project URL with the Zenodo DOI of the final report
#    in column T (URL) for every row that still references it.
#  - Strip the "D:/Ilich/" external-drive prefix from a few Server_Location
#    paths in column N, replacing it with a relative/NFWF-server style path.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldUrl = "https://www.marine.usf.edu/scamp/"
$newUrl = "https://zenodo.org/doi/10.5281/zenodo.8381009"

# Find the extent of the used range so we know how many rows to scan.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column T is the 20th column (URL column).
$urlCol = 20

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $urlCol)
    if ($cell.Value2 -eq $oldUrl) {
        $cell.Value2 = $newUrl
    }
}

# Column N is the 14th column (Server_Location).
$pathCol = 14
$oldPrefix = "D:/Ilich/"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $pathCol)
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().StartsWith($oldPrefix)) {
        $cell.Value2 = $val.ToString().Substring($oldPrefix.Length)
    }
}
